# Login.docx UCD edit:
#  1. Fill in the empty "Last Revision Date" value cell with "24/04/2017"
#     (dropping the cell's centered alignment) and add the transient
#     "_GoBack" bookmark right after the new run.
#  2. Collapse the run-per-word "In the activity 5 ... Not have this user
#     id ..." bullet into a single run with one consolidated w:t.
#  3. Collapse the run-per-word "In the activity 5 ... Incorrect password
#     ..." bullet into a single run with one consolidated w:t.
#  4. Remove the stray "_GoBack" bookmark that used to sit in the trailing
#     empty paragraph at the end of the document body.

$d = $word.ActiveDocument

# --- 1. Date Created / Last Revision Date row --------------------------
$table = $d.Tables.Item(1)
$dateCell = $table.Cell(4, 5)
$dateXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00EA4117" w:rsidRPr="00014107" w:rsidRDefault="00EA4117" w:rsidP="00014107"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>24/04/2017</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$dateCell.Range.InsertXML($dateXml)

# --- 2 & 3. Merge the fragmented runs in the two Alternative Flow bullets ---
$mergedNotSame = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BB77B9" w:rsidRDefault="00BB77B9" w:rsidP="00BB77B9"><w:pPr><w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>In the activity 5 of Normal flow, if the user id is not the same the data in the database system provide &#8220;Not have this user id in the database please input it again.&#8221; And back to the activity 2 of Normal flow.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$mergedWrongPass = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BB77B9" w:rsidRPr="00BB77B9" w:rsidRDefault="00BB77B9" w:rsidP="00BB77B9"><w:pPr><w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">In the activity 5 of Normal flow, if the user id is the same the data in the database but the password is not the same system provide &#8220;Incorrect password please input it again.&#8221; And back to the activity 2 of Normal flow.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# NOTE: after the structural InsertXML above, the cached $d.Paragraphs
# collection can resolve every index to the first paragraph; re-derive the
# paragraph collection from $d.Content (fresh) before walking it again.
$paraCount = $d.Content.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Content.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -like "In the act*not the same the data in the database*") {
        $para.Range.InsertXML($mergedNotSame)
        break
    }
}

$paraCount = $d.Content.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Content.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -like "In the activity 5 of Normal flow, if the user id is the same*Incorrect password*") {
        $para.Range.InsertXML($mergedWrongPass)
        break
    }
}

# --- 4. Drop the leftover "_GoBack" bookmark at the end of the document ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
